$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.759.50'
$ws.Range("E2").Value = '  +3.57%  '

$ws.Range("D3").Value = '3.129.18'
$ws.Range("E3").Value = '  +2.32%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.12'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.29'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.21%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '3.119.67'
$ws.Range("E8").Value = '  +2.35%  '

$ws.Range("E9").Value = '  +2.58%  '

$ws.Range("E10").Value = '  +19.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.70'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.83%  '

$ws.Range("E12").Value = '  +0.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +8.07%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.12'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.85%  '

$ws.Range("E15").Value = '  -0.26%  '

$ws.Range("D16").Value = '3.647.99'
$ws.Range("E16").Value = '  +2.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.21'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("D18").Value = '63.662.44'
$ws.Range("E18").Value = '  +3.40%  '

$ws.Range("D19").Value = '3.126.10'
$ws.Range("E19").Value = '  +2.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.31'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.26'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.735'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("E23").Value = '  +2.74%  '

$ws.Range("E24").Value = '  -1.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.49'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.98%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.69'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +8.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.72'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.83'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.08'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.53%  '

$ws.Range("E33").Value = '  +1.54%  '

$ws.Range("D34").Value = '0.0₃0870'
$ws.Range("E34").Value = '  +10.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.40'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +11.75%  '

$ws.Range("E36").Value = '  +2.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.43'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +16.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.14'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.87'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '453.34'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +7.57%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.74'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0373'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.90%  '

$ws.Range("D43").Value = '2.924.26'
$ws.Range("E43").Value = '  +5.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.276'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.49%  '

$ws.Range("E45").Value = '  +2.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.17'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '127.29'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.44%  '

$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.79'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.38%  '

$ws.Range("E50").Value = '  +0.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.67'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.12%  '
